$d = $word.ActiveDocument

$d.Content.Find.Execute(
    "Move forward 3 spaces.  If one of these spaces is occupied by an enemy, stop there and push the enemy back 1 space.  Deal 6 damage, take 1 damage if an enemy is contacted.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Move forward exactly 3 spaces in a straight line until contacting an enemy or border.  If contacting an enemy, stop there and push the enemy back 1 space.  Deal 6 damage, take 1 damage if an enemy is contacted.",
    2)

$d.Content.Find.Execute(
    "in any direction.  Enemies that occupy this space are knocked back 1 space and take 2 damage.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "in any direction.  Enemies that occupy this space are knocked back 1 space and take 2 damage.",
    2)

$d.Content.Find.Execute(
    "Enemies in the blast take 6 damage, you take 2, allies take 0.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Enemies in the blast take 6 damage, you take 2, allies take 0.",
    2)

$d.Content.Find.Execute(
    "Heal for 5 damage, draw 1 card.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Heal for 5 damage, draw 1 card.",
    2)
